$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 104 and 105: swap all data (columns B..AB), keep column A (row id) unchanged ---
$ws.Range("B104").Value = 7127370
$ws.Range("C104").Value = "Australia ALeague"
$ws.Range("D104").Value = 45340.125
$ws.Range("E104").Value = "Macarthur FC"
$ws.Range("F104").Value = "Wellington Phoenix"
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 2
$ws.Range("I104").Value = "A"
$ws.Range("J104").Value = 2.4
$ws.Range("K104").Value = 3.75
$ws.Range("L104").Value = 2.625
$ws.Range("M104").Value = 2.375
$ws.Range("N104").Value = 3.8
$ws.Range("O104").Value = 2.75
$ws.Range("P104").Value = 0
$ws.Range("Q104").Value = 1.8
$ws.Range("R104").Value = 2.05
$ws.Range("S104").Value = 3
$ws.Range("T104").Value = 1.9
$ws.Range("U104").Value = 1.95
$ws.Range("V104").Value = -1
$ws.Range("W104").Value = -1
$ws.Range("X104").Value = 1.75
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 1.05
$ws.Range("AA104").Value = 0
$ws.Range("AB104").Value = 0

$ws.Range("B105").Value = 7127374
$ws.Range("C105").Value = "Australia ALeague"
$ws.Range("D105").Value = 45340.125
$ws.Range("E105").Value = "Central Coast Mariners"
$ws.Range("F105").Value = "Western Sydney Wanderers"
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = "H"
$ws.Range("J105").Value = 1.909
$ws.Range("K105").Value = 3.75
$ws.Range("L105").Value = 3.6
$ws.Range("M105").Value = 2.15
$ws.Range("N105").Value = 3.6
$ws.Range("O105").Value = 3.25
$ws.Range("P105").Value = -0.25
$ws.Range("Q105").Value = 1.86
$ws.Range("R105").Value = 2.04
$ws.Range("S105").Value = 2.75
$ws.Range("T105").Value = 1.975
$ws.Range("U105").Value = 1.875
$ws.Range("V105").Value = 1.15
$ws.Range("W105").Value = -1
$ws.Range("X105").Value = -1
$ws.Range("Y105").Value = 0.8600000000000001
$ws.Range("Z105").Value = -1
$ws.Range("AA105").Value = -1
$ws.Range("AB105").Value = 0.875

# --- Rows 112 and 113: swap all data (columns B..AB), keep column A (row id) unchanged ---
$ws.Range("B112").Value = 7127376
$ws.Range("C112").Value = "Australia ALeague"
$ws.Range("D112").Value = 45347.125
$ws.Range("E112").Value = "Newcastle Jets"
$ws.Range("F112").Value = "Macarthur FC"
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = 2
$ws.Range("I112").Value = "D"
$ws.Range("J112").Value = 1.95
$ws.Range("K112").Value = 4
$ws.Range("L112").Value = 3.4
$ws.Range("M112").Value = 1.909
$ws.Range("N112").Value = 4.2
$ws.Range("O112").Value = 3.6
$ws.Range("P112").Value = -0.5
$ws.Range("Q112").Value = 1.89
$ws.Range("R112").Value = 2.01
$ws.Range("S112").Value = 3.5
$ws.Range("T112").Value = 1.95
$ws.Range("U112").Value = 1.9
$ws.Range("V112").Value = -1
$ws.Range("W112").Value = 3.2
$ws.Range("X112").Value = -1
$ws.Range("Y112").Value = -1
$ws.Range("Z112").Value = 1.01
$ws.Range("AA112").Value = 0.95
$ws.Range("AB112").Value = -1

$ws.Range("B113").Value = 7127379
$ws.Range("C113").Value = "Australia ALeague"
$ws.Range("D113").Value = 45347.125
$ws.Range("E113").Value = "Melbourne Victory"
$ws.Range("F113").Value = "Central Coast Mariners"
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 1
$ws.Range("I113").Value = "A"
$ws.Range("J113").Value = 1.95
$ws.Range("K113").Value = 3.6
$ws.Range("L113").Value = 3.8
$ws.Range("M113").Value = 1.909
$ws.Range("N113").Value = 3.6
$ws.Range("O113").Value = 4
$ws.Range("P113").Value = -0.5
$ws.Range("Q113").Value = 1.9
$ws.Range("R113").Value = 1.95
$ws.Range("S113").Value = 2.75
$ws.Range("T113").Value = 1.925
$ws.Range("U113").Value = 1.925
$ws.Range("V113").Value = -1
$ws.Range("W113").Value = -1
$ws.Range("X113").Value = 3
$ws.Range("Y113").Value = -1
$ws.Range("Z113").Value = 0.95
$ws.Range("AA113").Value = -1
$ws.Range("AB113").Value = 0.925

# --- Rows 124 and 125: swap all data (columns B..AB), keep column A (row id) unchanged ---
$ws.Range("B124").Value = 7127388
$ws.Range("C124").Value = "Australia ALeague"
$ws.Range("D124").Value = 45361.125
$ws.Range("E124").Value = "Sydney FC"
$ws.Range("F124").Value = "Brisbane Roar"
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 1
$ws.Range("I124").Value = "D"
$ws.Range("J124").Value = 1.5
$ws.Range("K124").Value = 5
$ws.Range("L124").Value = 5
$ws.Range("M124").Value = 1.533
$ws.Range("N124").Value = 5.25
$ws.Range("O124").Value = 5
$ws.Range("P124").Value = -1
$ws.Range("Q124").Value = 1.8
$ws.Range("R124").Value = 2.05
$ws.Range("S124").Value = 3.5
$ws.Range("T124").Value = 1.925
$ws.Range("U124").Value = 1.925
$ws.Range("V124").Value = -1
$ws.Range("W124").Value = 4.25
$ws.Range("X124").Value = -1
$ws.Range("Y124").Value = -1
$ws.Range("Z124").Value = 1.05
$ws.Range("AA124").Value = -1
$ws.Range("AB124").Value = 0.925

$ws.Range("B125").Value = 7128012
$ws.Range("C125").Value = "Australia ALeague"
$ws.Range("D125").Value = 45361.125
$ws.Range("E125").Value = "Macarthur FC"
$ws.Range("F125").Value = "Central Coast Mariners"
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 3
$ws.Range("I125").Value = "A"
$ws.Range("J125").Value = 2.4
$ws.Range("K125").Value = 3.5
$ws.Range("L125").Value = 2.75
$ws.Range("M125").Value = 3.4
$ws.Range("N125").Value = 3.75
$ws.Range("O125").Value = 2.05
$ws.Range("P125").Value = 0.25
$ws.Range("Q125").Value = 2.025
$ws.Range("R125").Value = 1.825
$ws.Range("S125").Value = 3
$ws.Range("T125").Value = 2.05
$ws.Range("U125").Value = 1.8
$ws.Range("V125").Value = -1
$ws.Range("W125").Value = -1
$ws.Range("X125").Value = 1.05
$ws.Range("Y125").Value = -1
$ws.Range("Z125").Value = 0.825
$ws.Range("AA125").Value = 0
$ws.Range("AB125").Value = 0

# --- Row 167: update odds-related cells only ---
$ws.Range("M167").Value = 2.05
$ws.Range("N167").Value = 3.8
$ws.Range("O167").Value = 3.3
$ws.Range("Q167").Value = 2.07
$ws.Range("R167").Value = 1.83
$ws.Range("T167").Value = 1.925
$ws.Range("U167").Value = 1.925

# --- Row 168: update odds-related cells only ---
$ws.Range("Q168").Value = 2.05
$ws.Range("R168").Value = 1.85
$ws.Range("T168").Value = 1.925
$ws.Range("U168").Value = 1.925
